# EI Variable Installments T1 scenarios
# Applies view/selection changes and adds the "O" (extra "Date"/variable-installment)
# column values on the Repayment schedule sheet, plus a column width tweak on Summary.

$wb = $excel.ActiveWorkbook

$wsNewLoan = $wb.Worksheets.Item("NewLoanInput")
$wsSummary = $wb.Worksheets.Item("Summary")
$wsRepay   = $wb.Worksheets.Item("Repayment schedule")
$wsTrans   = $wb.Worksheets.Item("Transactions")

# --- Repayment schedule: populate new column O (mirrors column N: style 16, value 0) ---
# Row 2 only carries the style (no value), rows 3-14 carry a numeric 0.
$wsRepay.Range("N2").Copy($wsRepay.Range("O2")) | Out-Null

for ($r = 3; $r -le 14; $r++) {
    $wsRepay.Range("N$r").Copy($wsRepay.Range("O$r")) | Out-Null
    $wsRepay.Range("O$r").Value = 0
}

# --- Summary: widen column B and drop the "best fit" auto-size flag ---
$wsSummary.Columns.Item(2).ColumnWidth = 6.88

# --- Selections on each sheet (also drives which sheet ends up "active") ---
$wsNewLoan.Range("B2").Select() | Out-Null
$wsSummary.Range("A8:XFD16").Select() | Out-Null
$wsRepay.Range("A15:XFD15").Select() | Out-Null
$wsTrans.Range("A2:XFD2").Select() | Out-Null
